# Generate Report for Handoff
#
# The localization status report moves from "In Translation" to
# "Ready for handoff" and the handoff timestamps are refreshed.
# Update the Overview roll-up sheet and each per-language sheet
# (zh-cn, de-de), then widen the Status / timestamp columns so the
# longer text still fits (AutoFit, as a real handoff-report generator
# would do right after rewriting the status text).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns + latest generate date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 21:05:46"

# zh-cn detail sheet: status + latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 21:05:41"

# de-de detail sheet: status + latest handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 21:05:46"

# Re-fit the columns that now hold the longer status/date text
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null
